# Applies the edits described in the diff for Saldo.xlsx / xl/worksheets/sheet.xml
# Net change:
#   - Insert a new row (004001621 / DANIELA / 217022.7) right above the LAILA row
#     (i.e. right after the PAULA row, as the new row 3).
#   - Remove the old DANIELA row (004001621 / DANIELA / 165512.75) that used to sit
#     right after the LAILA row.
#   - Remove the EULER row (004399832 / EULER / 16614.4) that used to sit right
#     before the VERA row.
#   - Update the VERA row's Saldo value from 15000 to 16671.2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand new row above row 3 (the LAILA row) and fill it with the new
#    DANIELA account/name/balance.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "004001621"
$ws.Range("B3").Value = "DANIELA"
$ws.Range("C3").Value = 217022.7

# After the insert above, the sheet looks like:
#   1 Conta/Nome/Saldo
#   2 PAULA
#   3 DANIELA  217022.7   (just inserted)
#   4 LAILA
#   5 DANIELA  165512.75  <- old row to remove
#   6 MARIO
#   7 RICARDO
#   8 THOMAS
#   9 EULER    16614.4    <- row to remove
#  10 VERA     15000      <- value to update

# 2) Remove the old DANIELA (165512.75) row, now at row 5.
$ws.Rows.Item(5).Delete()

# After removal, the sheet looks like:
#   1 Conta/Nome/Saldo
#   2 PAULA
#   3 DANIELA  217022.7
#   4 LAILA
#   5 MARIO
#   6 RICARDO
#   7 THOMAS
#   8 EULER    16614.4   <- row to remove
#   9 VERA     15000     <- value to update

# 3) Remove the EULER row, now at row 8.
$ws.Rows.Item(8).Delete()

# After removal, VERA is now at row 8.
#   1 Conta/Nome/Saldo
#   2 PAULA
#   3 DANIELA  217022.7
#   4 LAILA
#   5 MARIO
#   6 RICARDO
#   7 THOMAS
#   8 VERA     15000 -> 16671.2

# 4) Update VERA's Saldo value.
$ws.Range("C8").Value = 16671.2

Write-Host "Edit applied."
